$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value  = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G6").Value  = "Tumakuru (Tumkur)"
$ws.Range("G7").Value  = "Davangere"
$ws.Range("G8").Value  = "Dharwad"
$ws.Range("G14").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G16").Value = "Shivamogga (Shimoga)"
$ws.Range("G20").Value = "Davangere"
$ws.Range("G21").Value = "Dharwad"
$ws.Range("G23").Value = "Vijayapura (Bijapur)"
$ws.Range("G30").Value = "Vijayapura (Bijapur)"
$ws.Range("G31").Value = "Shivamogga (Shimoga)"
$ws.Range("G32").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G33").Value = "Bagalkot"
$ws.Range("G34").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G39").Value = "Vijayapura (Bijapur)"
$ws.Range("G41").Value = "Shivamogga (Shimoga)"
$ws.Range("G46").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G47").Value = "Vijayapura (Bijapur)"
$ws.Range("G50").Value = "Dharwad"
